# New Submission Synced: 2026-02-09 20:40:31
# A new Google-Forms-style response row is appended to the "JSS 3B" sheet,
# whose used range currently ends at row 8 (A1:D8 = header + 7 responses).
# The sheet dimension grows to A1:D9 and the new response lands on row 9.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3B")

$row = 9

# Timestamp (text, same free-form string format used by every other row)
$ws.Cells.Item($row, 1).Value = "2026-02-09 20:40:31"

# Full Name (text)
$ws.Cells.Item($row, 2).Value = "musa peter"

# Admission No - force text storage (matches the sheet's existing column C
# behaviour, which keeps admission numbers as text even when numeric-looking,
# e.g. "10", "36") so the numeric-looking "33" isn't auto-converted to a number.
$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = "33"

# AI Score (numeric)
$ws.Cells.Item($row, 4).Value = 9
